$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 4179.2
$ws.Range("I62").Value = 4322.294
$ws.Range("J62").Value = 3875.125
$ws.Range("K62").Value = 4322.294
$ws.Range("L62").Value = 3875.125
$ws.Range("M62").Value = -3698.294
$ws.Range("N62").Value = -5123.125
# Row 65
$ws.Range("H65").Value = 4179.2
$ws.Range("I65").Value = 4322.294
$ws.Range("J65").Value = 3875.125
$ws.Range("K65").Value = 21611.47
$ws.Range("L65").Value = 19375.625
$ws.Range("M65").Value = -18491.47
$ws.Range("N65").Value = -25615.625
# Row 92
$ws.Range("H92").Value = 29384.771
$ws.Range("I92").Value = 35386.93
$ws.Range("J92").Value = 374.33334
$ws.Range("K92").Value = 35386.93
$ws.Range("L92").Value = 374.33334
$ws.Range("M92").Value = -34138.93
$ws.Range("N92").Value = -2870.33334
# Row 98
$ws.Range("H98").Value = 1563.2059
$ws.Range("I98").Value = 919.26666
$ws.Range("J98").Value = 6392.75
$ws.Range("K98").Value = 919.26666
$ws.Range("L98").Value = 6392.75
$ws.Range("M98").Value = 578.73334
$ws.Range("N98").Value = -9388.75
# Row 122
$ws.Range("H122").Value = 1563.2059
$ws.Range("I122").Value = 919.26666
$ws.Range("J122").Value = 6392.75
$ws.Range("K122").Value = 2757.79998
$ws.Range("L122").Value = 19178.25
$ws.Range("M122").Value = -307.7999799999998
$ws.Range("N122").Value = -24078.25
# Row 137
$ws.Range("H137").Value = 2246.2222
$ws.Range("I137").Value = 1778.3636
$ws.Range("J137").Value = 2567.875
$ws.Range("K137").Value = 5335.0908
$ws.Range("L137").Value = 7703.625
$ws.Range("M137").Value = -2785.0908
$ws.Range("N137").Value = -12803.625
# Row 138
$ws.Range("H138").Value = 4421.9033
$ws.Range("I138").Value = 4854.5
$ws.Range("J138").Value = 4148.684
$ws.Range("K138").Value = 14563.5
$ws.Range("L138").Value = 12446.052
$ws.Range("M138").Value = -9423.5
$ws.Range("N138").Value = -22726.052

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 45
$ws.Range("H45").Value = 3827.5
$ws.Range("I45").Value = 3484
$ws.Range("J45").Value = 5316
$ws.Range("K45").Value = 3484
$ws.Range("L45").Value = 5316
$ws.Range("M45").Value = -3107
$ws.Range("N45").Value = -6070
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
# Row 61
$ws.Range("H61").Value = 6787.9287
$ws.Range("I61").Value = 6787.9287
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 6787.9287
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -6575.9287
$ws.Range("N61").ClearContents()
# Row 136
$ws.Range("H136").Value = 6787.9287
$ws.Range("I136").Value = 6787.9287
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 20363.7861
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -17813.7861
$ws.Range("N136").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1410.7142
$ws.Range("I16").Value = 1479.1666
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1479.1666
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1192.1666
$ws.Range("N16").Value = -1574
# Row 31
$ws.Range("H31").Value = 2031.0769
$ws.Range("I31").Value = 1995.6364
$ws.Range("J31").Value = 2057.0667
$ws.Range("K31").Value = 1995.6364
$ws.Range("L31").Value = 2057.0667
$ws.Range("M31").Value = -1700.6364
$ws.Range("N31").Value = -2647.0667
# Row 34
$ws.Range("H34").Value = 2031.0769
$ws.Range("I34").Value = 1995.6364
$ws.Range("J34").Value = 2057.0667
$ws.Range("K34").Value = 1995.6364
$ws.Range("L34").Value = 2057.0667
$ws.Range("M34").Value = -1793.6364
$ws.Range("N34").Value = -2461.0667
# Row 113
$ws.Range("H113").Value = 1410.7142
$ws.Range("I113").Value = 1479.1666
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 1479.1666
$ws.Range("L113").Value = 1000
$ws.Range("M113").Value = 690.8334
$ws.Range("N113").Value = -5340

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3827.6667
$ws.Range("I68").Value = 3419.8
$ws.Range("J68").Value = 3984.5386
$ws.Range("K68").Value = 10259.4
$ws.Range("L68").Value = 11953.6158
$ws.Range("M68").Value = -9448.400000000001
$ws.Range("N68").Value = -13575.6158
# Row 71
$ws.Range("H71").Value = 3827.6667
$ws.Range("I71").Value = 3419.8
$ws.Range("J71").Value = 3984.5386
$ws.Range("K71").Value = 30778.2
$ws.Range("L71").Value = 35860.8474
$ws.Range("M71").Value = -26722.2
$ws.Range("N71").Value = -43972.8474
# Row 107
$ws.Range("H107").Value = 3287.25
$ws.Range("J107").Value = 4133
$ws.Range("L107").Value = 12399
$ws.Range("N107").Value = -16239

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 5016.5
$ws.Range("I43").Value = 5016.5
$ws.Range("K43").Value = 5016.5
$ws.Range("M43").Value = -4865.5
# Row 46
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9844
# Row 80
$ws.Range("H80").Value = 12609.667
$ws.Range("I80").Value = 2980
$ws.Range("J80").Value = 17424.5
$ws.Range("K80").Value = 2980
$ws.Range("L80").Value = 17424.5
$ws.Range("M80").Value = -1982
$ws.Range("N80").Value = -19420.5
# Row 83
$ws.Range("H83").Value = 12609.667
$ws.Range("I83").Value = 2980
$ws.Range("J83").Value = 17424.5
$ws.Range("K83").Value = 14900
$ws.Range("L83").Value = 87122.5
$ws.Range("M83").Value = -9908
$ws.Range("N83").Value = -97106.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 159401.28
$ws.Range("I22").Value = 185718.33
$ws.Range("J22").Value = 1499
$ws.Range("K22").Value = 185718.33
$ws.Range("L22").Value = 1499
$ws.Range("M22").Value = -185423.33
$ws.Range("N22").Value = -2089
# Row 27
$ws.Range("H27").Value = 159401.28
$ws.Range("I27").Value = 185718.33
$ws.Range("J27").Value = 1499
$ws.Range("K27").Value = 185718.33
$ws.Range("L27").Value = 1499
$ws.Range("M27").Value = -185611.33
$ws.Range("N27").Value = -1713
# Row 55
$ws.Range("H55").Value = 916.6923
$ws.Range("I55").Value = 845.375
$ws.Range("J55").Value = 1030.8
$ws.Range("K55").Value = 845.375
$ws.Range("L55").Value = 1030.8
$ws.Range("M55").Value = -672.375
$ws.Range("N55").Value = -1376.8
# Row 61
$ws.Range("H61").Value = 6505.875
$ws.Range("I61").Value = 1661.75
$ws.Range("K61").Value = 1661.75
$ws.Range("M61").Value = -1459.75
# Row 113
$ws.Range("H113").Value = 6505.875
$ws.Range("I113").Value = 1661.75
$ws.Range("K113").Value = 1661.75
$ws.Range("M113").Value = 508.25
# Row 127
$ws.Range("H127").Value = 82250
$ws.Range("J127").Value = 82250
$ws.Range("L127").Value = 82250
$ws.Range("N127").Value = -92170
# Row 132
$ws.Range("H132").Value = 49368.92
$ws.Range("I132").Value = 62035.5
$ws.Range("J132").Value = 7147
$ws.Range("K132").Value = 186106.5
$ws.Range("L132").Value = 21441
$ws.Range("M132").Value = -183576.5
$ws.Range("N132").Value = -26501

